$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 360953.47
$ws.Range("J17").Value = 392271.2
$ws.Range("L17").Value = 1176813.6
$ws.Range("N17").Value = -1177149.6

$ws.Range("H69").Value = 4690.75
$ws.Range("I69").Value = 3763
$ws.Range("K69").Value = 11289
$ws.Range("M69").Value = -10415

$ws.Range("H72").Value = 4690.75
$ws.Range("I72").Value = 3763
$ws.Range("K72").Value = 33867
$ws.Range("M72").Value = -29499

$ws.Range("H87").Value = 15778.75
$ws.Range("I87").Value = 5555
$ws.Range("J87").Value = 16316.842
$ws.Range("K87").Value = 5555
$ws.Range("L87").Value = 16316.842
$ws.Range("M87").Value = -4307
$ws.Range("N87").Value = -18812.842

$ws.Range("H90").Value = 15778.75
$ws.Range("I90").Value = 5555
$ws.Range("J90").Value = 16316.842
$ws.Range("K90").Value = 16665
$ws.Range("L90").Value = 48950.526
$ws.Range("M90").Value = -10425
$ws.Range("N90").Value = -61430.526

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3892.5
$ws.Range("I61").Value = 2755.7334
$ws.Range("J61").Value = 7302.8
$ws.Range("K61").Value = 2755.7334
$ws.Range("L61").Value = 7302.8
$ws.Range("M61").Value = -2543.7334
$ws.Range("N61").Value = -7726.8

$ws.Range("H74").Value = 1802.3667
$ws.Range("I74").Value = 1131.8148
$ws.Range("J74").Value = 7837.3335
$ws.Range("K74").Value = 1131.8148
$ws.Range("L74").Value = 7837.3335
$ws.Range("M74").Value = -257.8148000000001
$ws.Range("N74").Value = -9585.333500000001

$ws.Range("H77").Value = 1802.3667
$ws.Range("I77").Value = 1131.8148
$ws.Range("J77").Value = 7837.3335
$ws.Range("K77").Value = 5659.074000000001
$ws.Range("L77").Value = 39186.6675
$ws.Range("M77").Value = -1291.074000000001
$ws.Range("N77").Value = -47922.6675

$ws.Range("H119").Value = 22449
$ws.Range("J119").Value = 22449
$ws.Range("L119").Value = 22449
$ws.Range("N119").Value = -32125

$ws.Range("H122").Value = 666.6129
$ws.Range("I122").Value = 522.85
$ws.Range("J122").Value = 928
$ws.Range("K122").Value = 1568.55
$ws.Range("L122").Value = 2784
$ws.Range("M122").Value = 881.4499999999998
$ws.Range("N122").Value = -7684

$ws.Range("H136").Value = 3892.5
$ws.Range("I136").Value = 2755.7334
$ws.Range("J136").Value = 7302.8
$ws.Range("K136").Value = 8267.200199999999
$ws.Range("L136").Value = 21908.4
$ws.Range("M136").Value = -5717.200199999999
$ws.Range("N136").Value = -27008.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1527.7273
$ws.Range("I107").Value = 1382.2963
$ws.Range("J107").Value = 2182.1667
$ws.Range("K107").Value = 1382.2963
$ws.Range("L107").Value = 2182.1667
$ws.Range("M107").Value = 537.7037
$ws.Range("N107").Value = -6022.1667

$ws.Range("H134").Value = 857.11475
$ws.Range("I134").Value = 854.3051
$ws.Range("J134").Value = 940
$ws.Range("K134").Value = 2562.9153
$ws.Range("L134").Value = 2820
$ws.Range("M134").Value = -27.91530000000012
$ws.Range("N134").Value = -7890

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1864.3715
$ws.Range("I58").Value = 1878.1786
$ws.Range("K58").Value = 1878.1786
$ws.Range("M58").Value = -1675.1786

$ws.Range("H99").Value = 2312.5
$ws.Range("I99").Value = 1655.5
$ws.Range("J99").Value = 2641
$ws.Range("K99").Value = 1655.5
$ws.Range("L99").Value = 2641
$ws.Range("M99").Value = -157.5
$ws.Range("N99").Value = -5637

$ws.Range("H126").Value = 2312.5
$ws.Range("I126").Value = 1655.5
$ws.Range("J126").Value = 2641
$ws.Range("K126").Value = 4966.5
$ws.Range("L126").Value = 7923
$ws.Range("M126").Value = -2496.5
$ws.Range("N126").Value = -12863

$ws.Range("H132").Value = 1696.4412
$ws.Range("I132").Value = 1541.9524
$ws.Range("J132").Value = 1946
$ws.Range("K132").Value = 4625.857199999999
$ws.Range("L132").Value = 5838
$ws.Range("M132").Value = -2095.857199999999
$ws.Range("N132").Value = -10898

$ws.Range("H134").Value = 1890.6097
$ws.Range("I134").Value = 1414.92
$ws.Range("J134").Value = 2633.875
$ws.Range("K134").Value = 4244.76
$ws.Range("L134").Value = 7901.625
$ws.Range("M134").Value = -1709.76
$ws.Range("N134").Value = -12971.625

$ws.Range("H136").Value = 1864.3715
$ws.Range("I136").Value = 1878.1786
$ws.Range("K136").Value = 5634.5358
$ws.Range("M136").Value = -3084.5358

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 58.827587
$ws.Range("I12").Value = 80.69231000000001
$ws.Range("J12").Value = 41.0625
$ws.Range("K12").Value = 242.07693
$ws.Range("L12").Value = 123.1875
$ws.Range("M12").Value = -69.07693
$ws.Range("N12").Value = -469.1875

$ws.Range("H34").Value = 1285.2
$ws.Range("J34").Value = 1382.3334
$ws.Range("L34").Value = 4147.0002
$ws.Range("N34").Value = -4315.0002

$ws.Range("H101").Value = 6989.25
$ws.Range("J101").Value = 6989.25
$ws.Range("L101").Value = 20967.75
$ws.Range("N101").Value = -25835.75

$ws.Range("H115").Value = 2147.3
$ws.Range("I115").Value = 1558.1666
$ws.Range("J115").Value = 3031
$ws.Range("K115").Value = 4674.4998
$ws.Range("L115").Value = 9093
$ws.Range("M115").Value = -3499.4998
$ws.Range("N115").Value = -11443

$ws.Range("H116").Value = 5145.125
$ws.Range("I116").Value = 2025.8
$ws.Range("K116").Value = 6077.4
$ws.Range("M116").Value = -2635.4

$ws.Range("H117").Value = 549.2857
$ws.Range("I117").Value = 224.16667
$ws.Range("J117").Value = 2500
$ws.Range("K117").Value = 672.50001
$ws.Range("L117").Value = 7500
$ws.Range("M117").Value = 2769.49999
$ws.Range("N117").Value = -14384

$ws.Range("H121").Value = 2382267.2
$ws.Range("I121").Value = 500
$ws.Range("J121").Value = 2779228.5
$ws.Range("K121").Value = 1500
$ws.Range("L121").Value = 8337685.5
$ws.Range("M121").Value = -190
$ws.Range("N121").Value = -8340305.5

$ws.Range("H122").Value = 1270.9524
$ws.Range("I122").Value = 400
$ws.Range("J122").Value = 1924.1666
$ws.Range("K122").Value = 3600
$ws.Range("L122").Value = 17317.4994
$ws.Range("M122").Value = -1150
$ws.Range("N122").Value = -22217.4994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1842.1111
$ws.Range("I126").Value = 1672.375
$ws.Range("J126").Value = 3200
$ws.Range("K126").Value = 5017.125
$ws.Range("L126").Value = 9600
$ws.Range("M126").Value = -2547.125
$ws.Range("N126").Value = -14540

$ws.Range("H132").Value = 2640.45
$ws.Range("I132").Value = 2281.9375
$ws.Range("J132").Value = 4074.5
$ws.Range("K132").Value = 6845.8125
$ws.Range("L132").Value = 12223.5
$ws.Range("M132").Value = -4315.8125
$ws.Range("N132").Value = -17283.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5205.75
$ws.Range("I40").Value = 4808.625
$ws.Range("K40").Value = 4808.625
$ws.Range("M40").Value = -4672.625

$ws.Range("H68").Value = 2061.8115
$ws.Range("I68").Value = 975.63336
$ws.Range("J68").Value = 2897.3333
$ws.Range("K68").Value = 975.63336
$ws.Range("L68").Value = 2897.3333
$ws.Range("M68").Value = -226.63336
$ws.Range("N68").Value = -4395.3333

$ws.Range("H71").Value = 2061.8115
$ws.Range("I71").Value = 975.63336
$ws.Range("J71").Value = 2897.3333
$ws.Range("K71").Value = 4878.1668
$ws.Range("L71").Value = 14486.6665
$ws.Range("M71").Value = -1134.1668
$ws.Range("N71").Value = -21974.6665

$ws.Range("H136").Value = 3476.4211
$ws.Range("I136").Value = 3139.111
$ws.Range("J136").Value = 3780
$ws.Range("K136").Value = 9417.332999999999
$ws.Range("L136").Value = 11340
$ws.Range("M136").Value = -6867.332999999999
$ws.Range("N136").Value = -16440

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 209.4375
$ws.Range("I113").Value = 169.27272
$ws.Range("K113").Value = 507.81816
$ws.Range("M113").Value = 1662.18184

$ws.Range("H126").Value = 626447.7
$ws.Range("I126").Value = 1000995.3
$ws.Range("J126").Value = 2201.6667
$ws.Range("K126").Value = 3002985.9
$ws.Range("L126").Value = 6605.000100000001
$ws.Range("M126").Value = -3000515.9
$ws.Range("N126").Value = -11545.0001

$ws.Range("H132").Value = 1064.6608
$ws.Range("J132").Value = 1855.4445
$ws.Range("L132").Value = 5566.333500000001
$ws.Range("M132").Value = -209.7022999999999
$ws.Range("N132").Value = -10626.3335

$ws.Range("H136").Value = 971.76
$ws.Range("I136").Value = 920.0952
$ws.Range("J136").Value = 1243
$ws.Range("K136").Value = 2760.2856
$ws.Range("L136").Value = 3729
$ws.Range("M136").Value = -210.2856000000002
$ws.Range("N136").Value = -8829
